# Fruta / hortaliza, semanal
# Insert a new weekly price observation as row 195, shifting the
# existing rows 195-276 down to 196-277.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 195 (pushes old rows 195..276 down to 196..277)
$ws.Rows.Item(195).EntireRow.Insert()

# Populate the new row with the latest weekly observation
$ws.Range("A195").Value = 4
$ws.Range("B195").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C195").Value = "Los Lagos"
$ws.Range("D195").Value = [DateTime]"2022-06-14"
$ws.Range("E195").Value = 10
$ws.Range("F195").Value = 100112043
$ws.Range("G195").Value = "Pepino ensalada"
$ws.Range("H195").Value = "Sin especificar"
$ws.Range("I195").Value = "Primera"
$ws.Range("J195").Value = 400
$ws.Range("K195").Value = 23000
$ws.Range("L195").Value = 23000
$ws.Range("M195").Value = 23000
$ws.Range("N195").Value = "`$/caja 60 unidades"
$ws.Range("O195").Value = "Región de Arica y Parinacota"
$ws.Range("P195").Value = 383
$ws.Range("Q195").Value = 60
$ws.Range("R195").Value = "Hortaliza"
